# Updates the ShipmentTrackNum (column C) values in rows 2-22, and the
# mirrored PackageTrackNum (column D) values where applicable, to the new
# tracking numbers per the "26th july 2022" commit.
#
# The new values are written as text (matching the original shared-string
# / text cell type) by temporarily forcing a text number format so Excel
# does not silently coerce the numeric-looking strings into numbers, then
# restoring each cell's original style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = "320018728611"
    3  = "320018728622"
    4  = "320018728655"
    5  = "320018728677"
    6  = "320018728714"
    7  = "320018728736"
    8  = "320018728770"
    9  = "320018728791"
    10 = "320018728840"
    11 = "320018728861"
    12 = "320018728909"
    13 = "320018728931"
    14 = "320018728975"
    15 = "320018729011"
    16 = "320018729044"
    17 = "320018729066"
    18 = "320018729103"
    19 = "320018729136"
    20 = "320018729170"
    21 = "320018729191"
    22 = "320018729228"
}

# Rows where column D (PackageTrackNum) mirrors column C (ShipmentTrackNum)
$dAlsoRows = @(5, 6, 7, 13, 14, 15, 16, 17)

function Set-TextValue($cell, $value) {
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $originalStyle
}

foreach ($row in ($newValues.Keys | Sort-Object)) {
    $value = $newValues[$row]

    $cCell = $ws.Cells.Item($row, 3)
    Set-TextValue $cCell $value

    if ($dAlsoRows -contains $row) {
        $dCell = $ws.Cells.Item($row, 4)
        Set-TextValue $dCell $value
    }
}
